$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new admin's e-mail address to the list.
$newRow = 9
$email = "amrinder.24bcs10596@sst.scaler.com"

$cell = $ws.Cells.Item($newRow, 1)
$cell.Value = $email

# Turn it into a mailto: hyperlink, matching the styling used for the
# other e-mail rows above it (rows 6-8 use the "Hyperlink" cell style).
$ws.Hyperlinks.Add($cell, "mailto:" + $email)
$cell.Style = $ws.Cells.Item($newRow - 1, 1).Style

# Move the active selection, matching the author's final cursor position.
[void]$ws.Range("D13").Select()
